# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 440 and 441) for
# Terminal Hortofrutícola Agro Chillán - Limón, pushing the existing
# rows 440-458 down to 442-460.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 440 (old data shifts down).
$ws.Rows.Item(440).Insert()
$ws.Rows.Item(440).Insert()

# Shared/common values for the two new rows (identical to the rest of the
# subset: same market, region, product taxonomy, unit, and kg/unidad).
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad = "Sin especificar"
$unidad = '$/malla 16 kilos'
$origen = "Región de O'Higgins"
$kgUnidad = 16
$fecha = 44509

# Row 440: "1a amarillo"
$r = 440
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "1a amarillo"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 6000
$ws.Cells.Item($r, 15).Value = 6500
$ws.Cells.Item($r, 16).Value = 6250
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 391
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 441: "2a amarillo"
$r = 441
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "2a amarillo"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 5000
$ws.Cells.Item($r, 15).Value = 5500
$ws.Cells.Item($r, 16).Value = 5250
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 328
$ws.Cells.Item($r, 20).Value = $kgUnidad
